$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header format from H1 (bold, bordered, centered) onto I1 and J1 only
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)

# Header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:I79 and J2:J79
$iValues = @(8,9,9,9,7,9,9,9,9,9,9,9,9,8,9,9,9,8,9,9,8,8,8,7,10,8,7,9,8,9,8,7,9,9,7,7,7,8,7,8,7,8,8,8,9,8,7,7,7,7,7,7,8,6,9,8,6,5,6,7,7,8,8,8,4,8,7,5,7,7,7,8,6,7,7,7,6,6)
$jValues = @(8,9,9,9,7,9,9,9,9,9,9,9,9,9,9,9,9,9,9,9,8,8,8,7,10,9,8,9,8,9,8,7,9,9,7,7,7,8,7,8,7,8,9,8,9,9,8,7,7,7,7,7,8,6,9,8,7,5,6,7,7,8,9,8,5,9,7,5,7,8,8,8,6,7,7,7,6,6)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}

